# Natural_Gas_Steo.xlsx -- update DataFeed values for the last two rows of
# 2022 (Nov/Dec already had correct historical data; Aug-Dec get refreshed)
# and all of 2023 (Jan-Dec), per the "Updated data feeds except mine netback
# and interim" commit. Also mirrors the formatting clean-up Excel performed
# on those same cells (the italic "preliminary" style is retired in favour
# of the normal bold numeric style / the plain default style), and restores
# the on-open selection/scroll position to where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Refresh the C-column (NGHHUUS) values for rows 177-193 (Aug 2022
#    through Dec 2023).
# ---------------------------------------------------------------------
$newValues = @{
    177 = 8.8000000000000007
    178 = 8.8570799999999998
    179 = 8.9706139999999994
    180 = 9.0125489999999999
    181 = 9.0938009999999991
    182 = 9.1019939999999995
    183 = 8.2639960000000006
    184 = 7.0404159999999996
    185 = 5.2307110000000003
    186 = 5.1922560000000004
    187 = 5.254505
    188 = 5.2992080000000001
    189 = 5.3116490000000001
    190 = 5.2273889999999996
    191 = 5.2520020000000001
    192 = 5.3757580000000003
    193 = 5.528791
}

foreach ($row in $newValues.Keys) {
    $ws.Range("C$row").Value = $newValues[$row]
}

# ---------------------------------------------------------------------
# 2. These rows had been using an italicized "preliminary estimate"
#    number style. With the data now final, that style is retired:
#    row 177 reverts to the plain default format (like the rest of the
#    un-styled C column), and rows 178-193 pick up the normal bold
#    right-aligned numeric format already used by C122:C176.
# ---------------------------------------------------------------------
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C177").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C122").Copy() | Out-Null
$ws.Range("C178:C193").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Restore the view: scrolled down to the bottom of the data with
#    E184 as the active (selected) cell.
# ---------------------------------------------------------------------
$ws.Range("A90").Select() | Out-Null
$ws.Range("E184").Select() | Out-Null
